$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking row): Right count and Wrong marking value
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total row): Right total, Wrong total, and the "X / Y" summary text
$ws.Range("B12").Value = 84
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "82 / 112"
